# global contact auto-save TC added
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- bookViews / window settings (best effort) ---
$excel.ActiveWindow.Height = 6350

# --- Update existing rows 13 & 14 (feature file rename + text updates) ---
$ws.Range("B13").Value = "globalContacts.feature"
$ws.Range("C13").Value = "SETUP: Launch Browser and go to application"

$ws.Range("B14").Value = "globalContacts.feature"
$ws.Range("C14").Value = "User verify user is on the Global Contact Creation page"
$ws.Range("E14").Value = "Yes"

# --- New rows 15-21: Global Contact test cases ---
$ws.Range("A15").Value = "TC_013"
$ws.Range("B15").Value = "globalContacts.feature"
$ws.Range("C15").Value = "Create and then Edit the same contact for Individual Global Contact"
$ws.Range("D15").Value = "No"
$ws.Range("E15").Value = "Yes"

$ws.Range("A16").Value = "TC_014"
$ws.Range("B16").Value = "globalContacts.feature"
$ws.Range("C16").Value = "Create the contact for Entity Global Contact"
$ws.Range("D16").Value = "No"
$ws.Range("E16").Value = "Yes"

$ws.Range("A17").Value = "TC_015"
$ws.Range("B17").Value = "globalContacts.feature"
$ws.Range("C17").Value = "Attempt to create a duplicate entity contact with the same EIN"
$ws.Range("D17").Value = "No"
$ws.Range("E17").Value = "Yes"

$ws.Range("A18").Value = "TC_016"
$ws.Range("B18").Value = "globalContacts.feature"
$ws.Range("C18").Value = "Verify Select & Proceed button is enabled after selecting a radio button"
$ws.Range("D18").Value = "No"
$ws.Range("E18").Value = "Yes"

$ws.Range("A19").Value = "TC_017"
$ws.Range("B19").Value = "globalContacts.feature"
$ws.Range("C19").Value = "Verify that the system trims leading and trailing spaces from text input fields"
$ws.Range("D19").Value = "No"
$ws.Range("E19").Value = "Yes"

$ws.Range("A20").Value = "TC_018"
$ws.Range("B20").Value = "globalContacts.feature"
$ws.Range("C20").Value = "Verify display result on entity name"
$ws.Range("D20").Value = "No"
$ws.Range("E20").Value = "Yes"

$ws.Range("A21").Value = "TC_019"
$ws.Range("B21").Value = "globalContacts.feature"
$ws.Range("C21").Value = "verify user authorization for `"View Only`" user"
$ws.Range("D21").Value = "No"
$ws.Range("E21").Value = "Yes"

# --- Selection moves to C14 ---
$ws.Range("C14").Select() | Out-Null
